$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("32:33").Insert()

$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 45054
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 100112021
$ws.Range("G32").Value = "Ají"
$ws.Range("H32").Value = "Inferno"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 140
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13500
$ws.Range("N32").Value = "`$/caja 15 kilos"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 900
$ws.Range("Q32").Value = 15
$ws.Range("R32").Value = "Hortaliza"

$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 45054
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100112021
$ws.Range("G33").Value = "Ají"
$ws.Range("H33").Value = "Inferno"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 160
$ws.Range("K33").Value = 9000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 9500
$ws.Range("N33").Value = "`$/caja 15 kilos"
$ws.Range("O33").Value = "Región de Arica y Parinacota"
$ws.Range("P33").Value = 633
$ws.Range("Q33").Value = 15
$ws.Range("R33").Value = "Hortaliza"
